$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row values (1:E1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 values
$ws.Range("B2").Value = 63.027085352698712
$ws.Range("C2").Value = 50.720884012982765
$ws.Range("D2").Value = 67.619524936529714
$ws.Range("E2").Value = 53.094941844320054

# Update row 3 values
$ws.Range("B3").Value = 64.919210033023205
$ws.Range("C3").Value = 45.78411335805194
$ws.Range("D3").Value = 76.406761011856247
$ws.Range("E3").Value = 51.476800363831472

# Update the selection to match the new range
$ws.Range("B1:E3").Select() | Out-Null
